# Generate Report for Handoff
# Updates the localization-status report: the source file was renamed
# (new GUID-named markdown source + new xliff content hash), the zh-cn
# handoff xliff was regenerated, and the de-de locale's "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" are reset
# because the new source hasn't been handed back yet.

$wb = $excel.ActiveWorkbook

$oldBase = "fb79a795-905e-4de4-91c2-ec580ac9116b"
$newBase = "e57dc085-a028-4d16-99de-e3959c44cce7"
$oldHash = "1e260ded4b1b56e70c10e3a6cc08507ffaaa9186"
$newHash = "934d7e9d0ab149d564528837a8a02e154c75e7f4"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

$wsO.Range("A2").Value = "$newBase.md"

# B2 carries a hyperlink; keep its target address but refresh the cell
# text / display text to the new file name.
$addrB2 = ""
foreach ($hl in $wsO.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$2") { $addrB2 = $hl.Address }
}
$wsO.Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), $addrB2, [Type]::Missing, [Type]::Missing, "e2e\$newBase.md") | Out-Null

$wsO.Range("G2").Value = "2016-08-19 00:57:19"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# A2 carries a hyperlink; keep address, refresh display text.
$addrA2 = ""
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") { $addrA2 = $hl.Address }
}
# I2 also carried its own hyperlink (rId3) which is being removed entirely.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $addrA2, [Type]::Missing, [Type]::Missing, "$newBase.md") | Out-Null

$wsZh.Range("G2").Value = "$newBase.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-19 00:57:14"

# Latest Target File / Latest Handback File are cleared (not yet handed
# back under the new source), Latest Handback DateTime resets to the
# "never" sentinel date.
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$addrA2de = ""
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") { $addrA2de = $hl.Address }
}
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $addrA2de, [Type]::Missing, [Type]::Missing, "$newBase.md") | Out-Null

$wsDe.Range("G2").Value = "$newBase.$newHash.de-de.xlf"
# Shares the same "Latest HO Xliff Generate Date" timestamp as the
# Overview sheet's G2 (both were regenerated together for de-de).
$wsDe.Range("H2").Value = "2016-08-19 00:57:19"

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
